$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-11 17:05:16.843000 to 2024-03-11 17:45:12.156000"
$ws.Range("B1").Style = "Normal"

$ws.Range("A2").Value = "Total time taken for the ride"
$ws.Range("B2").Value = 0.0277346412037037

$ws.Range("A3").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B3").Value = 20.89713222222223

$ws.Range("A4").Value = "Actual Watt-hours (Wh)"
$ws.Range("B4").Value = 1058.084981635556

$ws.Range("A5").Value = "Starting SoC (Ah)"
$ws.Range("B5").Value = 25.2

$ws.Range("A6").Value = "Ending SoC (Ah)"
$ws.Range("B6").Value = 3.427

$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = 63

$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = 8

$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("B9").Value = 24.58373603844754

$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("B10").Value = 43.04003996710556

$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = 55

$ws.Range("A12").Value = "Mode"
$ws.Range("B12").Value = "Custom mode`n56.98%`nEco mode`n43.02%"

$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("B13").Value = 5176.618119

$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("B14").Value = -1595.102987390285

$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B15").Value = 0.061052225

$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 0.005769735277205185

$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.338

$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 3.034

$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("B19").Value = 0.3040000000000003

$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("B20").Value = 40

$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("B21").Value = 47

$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 7

$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B23").Value = 66

$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B24").Value = 61

$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B25").Value = 62

$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("B26").Value = 68

$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("B27").Value = 98

$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B28").Value = 0

$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 47

$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = 40

$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B31").Value = 7

$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 53

$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.107548007777778

$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = [double]"1.284024308775944e-07"

$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 38

$ws.Range("A36").Value = "Idling time percentage"
$ws.Range("B36").Value = 1.320438906453413

$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = 4.863306676585457

$ws.Range("A38").Value = "Time spent in 10-20 km/h"
$ws.Range("B38").Value = 6.57429793565185

$ws.Range("A39").Value = "Time spent in 20-30 km/h"
$ws.Range("B39").Value = 12.66505486330668

$ws.Range("A40").Value = "Time spent in 30-40 km/h"
$ws.Range("B40").Value = 51.91091686814209

$ws.Range("A41").Value = "Time spent in 40-50 km/h"
$ws.Range("B41").Value = 9.842849172400966

$ws.Range("A42").Value = "Time spent in 50-60 km/h"
$ws.Range("B42").Value = 7.355402640877813

$ws.Range("A43").Value = "Time spent in 60-70 km/h"
$ws.Range("B43").Value = 4.379765668588432

$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 0.9717314487632509

$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0

$ws.Range("B2").NumberFormat = "[hh]:mm:ss"

Write-Host "done"